$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ODI Batting": the INNING_NUMBER (column B) cells that held no value
#    are now truly empty -- remove them so the cell no longer exists at all
#    (matches rows 2,3,4,6,11,12,13,14,15,16,17).
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$emptyInningRows = @(2, 3, 4, 6, 11, 12, 13, 14, 15, 16, 17)
foreach ($r in $emptyInningRows) {
    $batting.Range("B$r").ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" worksheet as the last tab (after
#    "ODI Bowling"), sheetId 4 / rId4.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------------
# 3. Header row: same bold / bordered / centered look as the other sheets.
#    Copy the formatting straight from an existing header cell so the new
#    sheet reuses the identical style.
# ---------------------------------------------------------------------------
$batting.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $ws.Cells.Item(1, $col + 1).Value = $headers[$col]
}

# ---------------------------------------------------------------------------
# 4. Data rows 2-19: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
#    PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH.
#    BATTING_POSITION (col B) holds real numbers; the rest are text, even
#    when the text looks numeric ("0", "1", "0.31%") so it must not be
#    reinterpreted as a number/percentage. Forcing a number format on every
#    cell -- including the blank ones -- also keeps those blank cells in
#    the sheet instead of having them silently dropped.
# ---------------------------------------------------------------------------
$ws.Range("A2:A19").NumberFormat = "@"
$ws.Range("B2:B19").NumberFormat = "General"
$ws.Range("C2:F19").NumberFormat = "@"

$data = @(
    @("4432", "", "", "", "", "NO"),
    @("4433", "", "", "", "", "NO"),
    @("4458", 10, "0", "0", "0.31%", "NO"),
    @("4459", 11, "", "", "", "NO"),
    @("4460", 11, "0", "0", "", "NO"),
    @("4472", 11, "0", "0", "", "NO"),
    @("4473", 10, "1", "0", "3.11%", "NO"),
    @("4476", "", "", "", "", "NO"),
    @("4564", 10, "", "", "", "NO"),
    @("4565", "", "", "", "", "NO"),
    @("4567", "", "", "", "", "NO"),
    @("4586", "", "", "", "", "NO"),
    @("4590", 10, "", "", "", "NO"),
    @("4634", "", "", "", "", "NO"),
    @("4638", "", "", "", "", "NO"),
    @("4686", "", "", "", "", ""),
    @("4688", "", "", "", "", ""),
    @("4690", "", "", "", "", "")
)

$rowIndex = 2
foreach ($record in $data) {
    for ($col = 0; $col -lt $record.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $record[$col]
    }
    $rowIndex++
}
